$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 6: Uur=5, project versie="1.4", Beschrijving=<rich text bold+plain run>
# ---------------------------------------------------------------------------

# A6 - plain number
$ws.Range("A6").Value = 5

# B6 - "1.4" stored as *text* (not coerced to a number), matching the style
# used by the other "project versie" cells (B2:B5). We first land the text
# via a helper formula cell (so Excel keeps it as a string) and paste only
# the value into B6, then copy the number/cell format from B5 onto B6.
$ws.Range("ZZ1").Formula = "=""1.4"""
$ws.Range("ZZ1").Copy()
$ws.Range("B6").PasteSpecial(-4163) | Out-Null   # xlPasteValues
$ws.Range("ZZ1").Clear()

$ws.Range("B5").Copy()
$ws.Range("B6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# C6 - rich text: bold lead-in sentence followed by the regular description,
# both using plain black Calibri 11pt.
$boldPart = "Officieel klaar met de user stories!! "
$restPart = "Dit uur heb ik ervoor gezorgd dat er particle effects worden gegeven aan de cells op de rand van de maze, en er word nu een animatie gespeeld wanneer een cell op beeld komt. Ik heb ook een maximum ingesteld zodat de X en Y niet groter kunnen zijn dan 250. het maken van de maze op deze grote is echter niet aan te raden. Ik kom hier later op terug in het README document. De user kan nu ook een nieuwe maze aanmaken waarneer die wilt."
$fullText = $boldPart + $restPart

$ws.Range("C6").Value = $fullText

# Start from the same format as the other Beschrijving cells (C2:C5).
$ws.Range("C5").Copy()
$ws.Range("C6").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# Make the cell's base font bold (this becomes the cell-level style).
[void]($ws.Range("C6").Font.Bold = $true)

# Explicitly format each run: bold lead-in, regular remainder.
$run1 = $ws.Range("C6").Characters(1, $boldPart.Length)
[void]($run1.Font.Bold = $true)
[void]($run1.Font.Name = "Calibri")
[void]($run1.Font.Size = 11)
[void]($run1.Font.Color = 0)

$run2 = $ws.Range("C6").Characters($boldPart.Length + 1, $restPart.Length)
[void]($run2.Font.Bold = $false)
[void]($run2.Font.Name = "Calibri")
[void]($run2.Font.Size = 11)
[void]($run2.Font.Color = 0)

# ---------------------------------------------------------------------------
# Sheet view: scrolled over to column B, selection now on C24
# ---------------------------------------------------------------------------
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("C24").Select() | Out-Null
